$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "OK"
$ws.Range("B5").Value = "apagada"
$ws.Range("C5").Value = "encendida"
$ws.Range("B6").Value = "apagado"
$ws.Range("C6").Value = "monitor roto"
$ws.Range("B17").Value = "d12d1212d12d"
